# Updated via Streamlit Approval System
# Appends 6 new pending-approval rows (rows 33-38) for
# "Western Interior Designers & Marine Contractors" to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 33
$endRow = 38

# --- Columns that are identical across all six new rows ---------------
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 1).Value  = "WGG 02"                                                      # A  EMPLOYEE_ID
    $ws.Cells.Item($r, 2).Value  = "Western Interior Designers & Marine Contractors"              # B  PARTY_NAME
    $ws.Cells.Item($r, 3).Value  = "16-01-2026"                                                   # C  DATE
    $ws.Cells.Item($r, 4).Value  = 286962                                                         # D  CORPORATE ID
    $ws.Cells.Item($r, 5).Value  = "Western Interior Designers & Marine Contractors"              # E  COMPANY NAME
    $ws.Cells.Item($r, 6).Value  = 34413429360                                                    # F  COMPANY ACCOUNT NO
    $ws.Cells.Item($r, 7).Value  = "NEFT"                                                         # G  TRANSACTION TYPE
    $ws.Cells.Item($r, 8).Value  = "SBIN0003229"                                                  # H  COMPANY IFSC
    $ws.Cells.Item($r, 9).Value  = "AAAFW8862C"                                                   # I  COMPANY PAN
    $ws.Cells.Item($r, 10).Value = "32AAAFW8862C1Z9"                                              # J  COMPANY GSTIN

    $ws.Cells.Item($r, 21).Value = "pending"                                                      # U  STATUS

    $ws.Cells.Item($r, 26).Value = "PAYMENT"                                                      # Z  CATEGORY
    $ws.Cells.Item($r, 27).Value = "Payments@westernidc.com"                                      # AA FROM_MAIL
    $ws.Cells.Item($r, 28).Value = "ESTIMATION NOT MATCHED"                                       # AB STATUS_MATCHED_ESTIMATION
    $ws.Cells.Item($r, 29).Value = 0                                                              # AC FIXED_AMOUNT
    $ws.Cells.Item($r, 30).Value = 0                                                              # AD BALANCE_AMOUNT
    $ws.Cells.Item($r, 31).Value = 0                                                              # AE ADJUSTMENT_AMOUNT
}

# --- Columns that differ per row: L (UNIQUE_ID), V (BASIC_AMOUNT), ----
# --- X (NARRATION), Y (PROJECT_NAME) ----------------------------------
$rowData = @(
    @("40e65b3d-445c-4745-9c7a-542174640e02", 126000, "Kolkata", "Kolkata RPA_UNIQUE_ID : 6d426d07-c34a-4ae3-aa78-aa9578353cd1"),
    @("0a40aecc-bc07-47f5-ae2b-8653728973c1", 66000,  "HO",      "Income tax payamet  Hijas Sir 2024-25 RPA_UNIQUE_ID : 74465b13-7088-4030-934a-09dc028a4422"),
    @("ed826686-f8ce-4b1f-8d60-643082983507", 66000,  "HO",      "Income tax payamet  Hisham Sir 2024-25 RPA_UNIQUE_ID : 835bc023-578f-4d70-b2dc-c7363cc9dbf0"),
    @("504a57a1-3288-44b8-b047-c6b21ea639df", 30000,  "HO",      "Income tax payamet Western 2024-25 RPA_UNIQUE_ID : 95b2ebbb-24b7-4fea-8752-7061975cc952"),
    @("bd9d67d4-1cf0-4f9c-ae69-d056eb7dda46", 38626,  "HO",      "SIDBI Due RPA_UNIQUE_ID : 3d4f589c-7892-4c16-a18d-cf60ec2b4dbe"),
    @("55326f55-6e97-4b53-b8e2-267e18a4b2ae", 1500,   "ho",      "Hisham sir advocate fee RPA_UNIQUE_ID : 68088f52-a26e-4145-8d9b-a409c34e325f")
)

$r = $startRow
foreach ($row in $rowData) {
    $ws.Cells.Item($r, 12).Value = $row[0]   # L  UNIQUE_ID
    $ws.Cells.Item($r, 22).Value = $row[1]   # V  BASIC_AMOUNT
    $ws.Cells.Item($r, 25).Value = $row[2]   # Y  PROJECT_NAME
    $ws.Cells.Item($r, 24).Value = $row[3]   # X  NARRATION
    $r++
}
